$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GoogleSearchTest")

# Update the test-result cells from FAIL to PASS
$ws.Range("F2").Value = "PASS"
$ws.Range("F3").Value = "PASS"

# Move the active selection to F3 to match the saved view state
$ws.Activate()
$ws.Range("F3").Select()
